$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "right" / "down" / "left" header columns (C/D/E) ---
# Before: C1=right, D1=down, E1=left
# After:  C1=down,  D1=left, E1=right
$ws.Range("C1").Value = "down"
$ws.Range("D1").Value = "left"
$ws.Range("E1").Value = "right"

# --- Move the data that lived under the old "right"/"left" columns to their
#     new homes, and fill in newly-populated cells from the updated map data ---

# Row 2: "right" value (2) moves from C2 to E2
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 2

# Row 3: "right" value (11) moves from C3 to E3
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 11

# Row 4: "left" value (12) moves from E4 to D4; "up" value (4) is newly added to B4
$ws.Range("E4").ClearContents()
$ws.Range("B4").Value = 4
$ws.Range("D4").Value = 12

# Row 5: new "down" value (3) in C5
$ws.Range("C5").Value = 3

# Row 12: new "left" value (2) in D12
$ws.Range("D12").Value = 2

# --- Add the documentation comments (mapping hotspots/regions to cells) ---
$c = $ws.Range("E2").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Hotspot : hExit1")

$c = $ws.Range("B3").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Hotspot : hExit1")

$c = $ws.Range("E3").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Region : region1")

$c = $ws.Range("B4").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Hotspot : hExit1")

$c = $ws.Range("D4").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Region : region1")

$c = $ws.Range("C5").AddComment()
$c.Text("Jack Mangano:" + [char]10 + "Hotspot : hExit1")

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("F10").Select()
